$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32: Automata for the People
$ws.Range("H32").Value = 1694.4
$ws.Range("J32").Value = 1490
$ws.Range("L32").Value = 1490
$ws.Range("N32").Value = -2142

# Row 113: Amaro Kart
$ws.Range("H113").Value = 2899.9092
$ws.Range("I113").Value = 2800
$ws.Range("J113").Value = 2957
$ws.Range("K113").Value = 2800
$ws.Range("L113").Value = 2957
$ws.Range("M113").Value = 454
$ws.Range("N113").Value = -9465

# Row 116: Growing Up
$ws.Range("H116").Value = 3166.6667
$ws.Range("J116").Value = 3500
$ws.Range("L116").Value = 3500
$ws.Range("N116").Value = -10384

# Row 125: Body over Mind
$ws.Range("H125").Value = 1669.2142
$ws.Range("I125").Value = 879
$ws.Range("J125").Value = 1884.7273
$ws.Range("K125").Value = 7911
$ws.Range("L125").Value = 16962.5457
$ws.Range("M125").Value = -5451
$ws.Range("N125").Value = -21882.5457

# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 2342.795
$ws.Range("I132").Value = 2141.3157
$ws.Range("J132").Value = 9999
$ws.Range("K132").Value = 6423.9471
$ws.Range("L132").Value = 29997
$ws.Range("M132").Value = -3893.9471
$ws.Range("N132").Value = -35057

# Row 138: All-night Crafting
$ws.Range("H138").Value = 2206.4666
$ws.Range("I138").Value = 1967.0416
$ws.Range("J138").Value = 3164.1667
$ws.Range("K138").Value = 5901.1248
$ws.Range("L138").Value = 9492.500100000001
$ws.Range("M138").Value = -761.1247999999996
$ws.Range("N138").Value = -19772.5001

$ws = $wb.Worksheets.Item("ARM")
# Row 41: Skillet Scandal
$ws.Range("H41").Value = 1389.25
$ws.Range("I41").Value = 1389.25
$ws.Range("K41").Value = 1389.25
$ws.Range("M41").Value = -975.25

# Row 45: Hollow Hallmarks
$ws.Range("H45").Value = 2679.8
$ws.Range("I45").Value = 2349.75
$ws.Range("K45").Value = 2349.75
$ws.Range("M45").Value = -1972.75

# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 2281.9768
$ws.Range("I61").Value = 1870.0303
$ws.Range("J61").Value = 3641.4
$ws.Range("K61").Value = 1870.0303
$ws.Range("L61").Value = 3641.4
$ws.Range("M61").Value = -1658.0303
$ws.Range("N61").Value = -4065.4

# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 1087.9166
$ws.Range("I74").Value = 745.2963
$ws.Range("K74").Value = 745.2963
$ws.Range("M74").Value = 128.7037

# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 1087.9166
$ws.Range("I77").Value = 745.2963
$ws.Range("K77").Value = 3726.4815
$ws.Range("M77").Value = 641.5185000000001

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 2281.9768
$ws.Range("I136").Value = 1870.0303
$ws.Range("J136").Value = 3641.4
$ws.Range("K136").Value = 5610.090899999999
$ws.Range("L136").Value = 10924.2
$ws.Range("M136").Value = -3060.090899999999
$ws.Range("N136").Value = -16024.2

$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin
$ws.Range("H86").Value = 2726.0908
$ws.Range("I86").Value = 2622.5
$ws.Range("J86").Value = 3002.3333
$ws.Range("K86").Value = 2622.5
$ws.Range("L86").Value = 3002.3333
$ws.Range("M86").Value = -1499.5
$ws.Range("N86").Value = -5248.3333

# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Range("H89").Value = 2726.0908
$ws.Range("I89").Value = 2622.5
$ws.Range("J89").Value = 3002.3333
$ws.Range("K89").Value = 13112.5
$ws.Range("L89").Value = 15011.6665
$ws.Range("M89").Value = -7496.5
$ws.Range("N89").Value = -26243.6665

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 2371.75
$ws.Range("I134").Value = 1952.2
$ws.Range("J134").Value = 3325.2727
$ws.Range("K134").Value = 5856.6
$ws.Range("L134").Value = 9975.8181
$ws.Range("M134").Value = -3321.6
$ws.Range("N134").Value = -15045.8181

$ws = $wb.Worksheets.Item("CRP")
# Row 80: The Long Armillae of the Law
$ws.Range("H80").Value = 12000
$ws.Range("I80").Value = 12000
$ws.Range("K80").Value = 12000
$ws.Range("M80").Value = -10877

# Row 83: Wooden Ambitions (L)
$ws.Range("H83").Value = 12000
$ws.Range("I83").Value = 12000
$ws.Range("K83").Value = 36000
$ws.Range("M83").Value = -30384

# Row 105: Zelkova, My Love
$ws.Range("H105").Value = 1699.5
$ws.Range("I105").Value = 400
$ws.Range("K105").Value = 400
$ws.Range("M105").Value = 1347

# Row 107: Built to Last
$ws.Range("H107").Value = 1320
$ws.Range("I107").Value = 600
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 600
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = 1320
$ws.Range("N107").Value = -5340

# Row 122: Timber of Tenkonto
$ws.Range("H122").Value = 1802.4517
$ws.Range("I122").Value = 1487.3334
$ws.Range("K122").Value = 4462.0002
$ws.Range("M122").Value = -2012.0002

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 5381.1035
$ws.Range("J134").Value = 5803.8
$ws.Range("L134").Value = 17411.4
$ws.Range("N134").Value = -22481.4

$ws = $wb.Worksheets.Item("CUL")
# Row 34: Fever Pitch
$ws.Range("H34").Value = 9804362
$ws.Range("I34").Value = 147
$ws.Range("J34").Value = 10638764
$ws.Range("K34").Value = 441
$ws.Range("L34").Value = 31916292
$ws.Range("M34").Value = -357
$ws.Range("N34").Value = -31916460

# Row 39: Bloody Good Tart, This
$ws.Range("H39").Value = 2035.5714
$ws.Range("I39").Value = 500
$ws.Range("J39").Value = 2454.3635
$ws.Range("K39").Value = 1500
$ws.Range("L39").Value = 7363.0905
$ws.Range("M39").Value = -1206
$ws.Range("N39").Value = -7951.0905

# Row 55: Pagan Pastries
$ws.Range("H55").Value = 1983.3334
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 1983.3334
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 5950.0002
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -6304.0002

# Row 141: Ocean Explosion
$ws.Range("H141").Value = 8865.888999999999
$ws.Range("I141").Value = 2506.6667
$ws.Range("J141").Value = 10137.733
$ws.Range("K141").Value = 7520.000100000001
$ws.Range("L141").Value = 30413.199
$ws.Range("M141").Value = -2340.000100000001
$ws.Range("N141").Value = -40773.199

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 50941480
$ws.Range("I80").Value = 101800760
$ws.Range("K80").Value = 101800760
$ws.Range("M80").Value = -101799762

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 50941480
$ws.Range("I83").Value = 101800760
$ws.Range("K83").Value = 509003800
$ws.Range("M83").Value = -508998808

# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 2276
$ws.Range("I122").Value = 1552
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 4656
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -2206
$ws.Range("N122").Value = -13900

# Row 132: On Board for Lar
$ws.Range("H132").Value = 3923.84
$ws.Range("I132").Value = 3475.647
$ws.Range("J132").Value = 4876.25
$ws.Range("K132").Value = 10426.941
$ws.Range("L132").Value = 14628.75
$ws.Range("M132").Value = -7896.940999999999
$ws.Range("N132").Value = -19688.75

$ws = $wb.Worksheets.Item("LTW")
# Row 5: These Boots Are Made for Wailing
$ws.Range("H5").Value = 62673.332
$ws.Range("J5").Value = 62673.332
$ws.Range("L5").Value = 62673.332
$ws.Range("N5").Value = -62899.332

# Row 43: Subordinate Clause
$ws.Range("H43").Value = 68601.39999999999
$ws.Range("J43").Value = 68601.39999999999
$ws.Range("L43").Value = 68601.39999999999
$ws.Range("N43").Value = -68987.39999999999

# Row 68: You Could Say It's a Moving Target
$ws.Range("H68").Value = 1500
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()

# Row 71: They Call It Bloody Mary (L)
$ws.Range("H71").Value = 1500
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()

# Row 82: Trainin' the Neck
$ws.Range("H82").Value = 45456044
$ws.Range("I82").Value = 50001450
$ws.Range("K82").Value = 50001450
$ws.Range("M82").Value = -50001089

# Row 85: Training Is Only Skintight (L)
$ws.Range("H85").Value = 45456044
$ws.Range("I85").Value = 50001450
$ws.Range("K85").Value = 50001450
$ws.Range("M85").Value = -50000202

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 1793.0233
$ws.Range("I132").Value = 1321.7097
$ws.Range("J132").Value = 3010.5833
$ws.Range("K132").Value = 3965.1291
$ws.Range("L132").Value = 9031.749899999999
$ws.Range("M132").Value = -1435.1291
$ws.Range("N132").Value = -14091.7499

# Row 133: The Perfect Accessory
$ws.Range("H133").Value = 34417.332
$ws.Range("J133").Value = 34417.332
$ws.Range("L133").Value = 34417.332
$ws.Range("N133").Value = -39477.332

# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 4386999
$ws.Range("I136").Value = 924.23334
$ws.Range("K136").Value = 2772.70002
$ws.Range("M136").Value = -222.7000200000002

$ws = $wb.Worksheets.Item("WVR")
# Row 56: Full Moon Fever
$ws.Range("H56").Value = 44500
$ws.Range("J56").Value = 44500
$ws.Range("L56").Value = 44500
$ws.Range("N56").Value = -45928

# Row 62: Pride Up in Smoke
$ws.Range("H62").Value = 34803.637
$ws.Range("I62").Value = 1900
$ws.Range("J62").Value = 38094
$ws.Range("K62").Value = 1900
$ws.Range("L62").Value = 38094
$ws.Range("M62").Value = -1276
$ws.Range("N62").Value = -39342

# Row 65: Desperate for Diversionaries (L)
$ws.Range("H65").Value = 34803.637
$ws.Range("I65").Value = 1900
$ws.Range("J65").Value = 38094
$ws.Range("K65").Value = 9500
$ws.Range("L65").Value = 190470
$ws.Range("M65").Value = -6380
$ws.Range("N65").Value = -196710

# Row 100: Of Great Import
$ws.Range("H100").Value = 1502.5625
$ws.Range("I100").Value = 1646.4546
$ws.Range("J100").Value = 1186
$ws.Range("K100").Value = 3292.9092
$ws.Range("L100").Value = 2372
$ws.Range("M100").Value = -2751.9092
$ws.Range("N100").Value = -3454

# Row 107: Flax Wax
$ws.Range("H107").Value = 574.3
$ws.Range("I107").Value = 709
$ws.Range("J107").Value = 372.25
$ws.Range("K107").Value = 2127
$ws.Range("L107").Value = 1116.75
$ws.Range("M107").Value = -207
$ws.Range("N107").Value = -4956.75

# Row 122: Heavy Armoire
$ws.Range("H122").Value = 1790.1
$ws.Range("I122").Value = 1758.6333
$ws.Range("J122").Value = 1884.5
$ws.Range("K122").Value = 5275.8999
$ws.Range("L122").Value = 5653.5
$ws.Range("M122").Value = -2825.8999
$ws.Range("N122").Value = -10553.5

# Row 126: A Polished Purchase
$ws.Range("H126").Value = 3207.4285
$ws.Range("I126").Value = 1300.6666
$ws.Range("J126").Value = 4637.5
$ws.Range("K126").Value = 3901.9998
$ws.Range("L126").Value = 13912.5
$ws.Range("M126").Value = -1431.9998
$ws.Range("N126").Value = -18852.5

# Row 130: Skill Cap
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 2102
$ws.Range("I136").Value = 1542.5428
$ws.Range("K136").Value = 4627.6284
$ws.Range("M136").Value = -2077.6284
